$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "19/9"
$ws.Range("A3").Value = "20/9"
$ws.Range("B2").Value = "D19CQDT03-B + D19CQDT04-B  "
$ws.Range("E2").Value = "D19CQDT01-B + D19CQDT02-B   "
